$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F26").Value = "Cell F26"
$ws.Range("H1").Value = "Cell H1"
$ws.Range("H2").Value = "Cell H2"
$ws.Range("H3").Value = "Cell H3"
